$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change header of column B from "Concept" to "Text"
$ws.Range("B1").Value = "Text"

# Replace the concept values in column B (rows 2-5) with "validity"
$ws.Range("B2").Value = "validity"
$ws.Range("B3").Value = "validity"
$ws.Range("B4").Value = "validity"
$ws.Range("B5").Value = "validity"

# Update active selection to B5
$ws.Range("B5").Select()
